# Renumber the section titles on slides 4-10 down by one (the deck used to
# skip "3", e.g. 1, 2, 4, 5, 6, 7, 8, 9, 10 -- this closes the gap so it
# reads 1, 2, 3, 4, 5, 6, 7, 8, 9).

$p = $ppt.ActivePresentation

# Slide 4: title run "4 " -> "3 "
$s4 = $p.Slides.Item(4)
$t4 = $s4.Shapes.Item(1).TextFrame.TextRange
$t4.Characters(1, 2).Text = "3 "

# Slide 5: title run "5 " -> "4 " (PowerPoint split this into a fresh "4"
# run followed by the pre-existing " " run when it was hand-edited)
$s5 = $p.Slides.Item(5)
$t5 = $s5.Shapes.Item(1).TextFrame.TextRange
$t5.Characters(1, 1).Text = "4"

# Slide 6: title run "6 " -> "5 "
$s6 = $p.Slides.Item(6)
$t6 = $s6.Shapes.Item(1).TextFrame.TextRange
$t6.Characters(1, 2).Text = "5 "

# Slide 7: title run "7 " -> "6 "
$s7 = $p.Slides.Item(7)
$t7 = $s7.Shapes.Item(1).TextFrame.TextRange
$t7.Characters(1, 2).Text = "6 "

# Slide 8: title run "8 " -> "7 "
$s8 = $p.Slides.Item(8)
$t8 = $s8.Shapes.Item(1).TextFrame.TextRange
$t8.Characters(1, 2).Text = "7 "

# Slide 9: title run "9 " -> "8 "
$s9 = $p.Slides.Item(9)
$t9 = $s9.Shapes.Item(1).TextFrame.TextRange
$t9.Characters(1, 2).Text = "8 "

# Slide 10: title run "10 " -> "9 "
$s10 = $p.Slides.Item(10)
$t10 = $s10.Shapes.Item(1).TextFrame.TextRange
$t10.Characters(1, 3).Text = "9 "
